$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8").Value = "Meilleure: 1095445.00`nMoyenne: 1094931.07`nÉcart type: 5997185.54"
$ws.Range("G9").Value = "Meilleure: 624319.00`nMoyenne: 620104.20`nÉcart type: 3396464.64"
$ws.Range("G10").Value = "Meilleure: 7772.00`nMoyenne: 7760.37`nÉcart type: 42505.32"
